$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the now-unused trailing rows 24-26 first --------------
# (Bibliografia:/A ser definida.../Requisitos:/LOM3110... content), so the
# sheet dimension shrinks from A1:C26 to A1:C23 and row numbers below stay
# put while we rewrite rows 13-23.
$ws.Range("A24:C26").EntireRow.Delete()

# --- Step 2: rewrite rows 13-23 -------------------------------------------
# The "Programa resumido" block moved up 3 rows, the Docentes responsaveis
# names now double up into the B/C columns of rows 13/15/18, and the old
# long "Elaboracao..." / "A disciplina consiste..." paragraphs are dropped.

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Rows.Item(13).RowHeight = 60

$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

$ws.Range("A14").Value = "Short syllabus:"
$ws.Rows.Item(14).RowHeight = 60

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C15").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Rows.Item(15).RowHeight = 120

$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Rows.Item(18).RowHeight = 60

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Range("C19").Value = "O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas."
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(19).RowHeight = 60

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Range("C20").Value = "A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina."
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(20).RowHeight = 60

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "não há"
$ws.Range("C21").Value = "não há"
$ws.Rows.Item(21).RowHeight = 120

$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

$ws.Range("A22").Value = "Requisitos:"
$ws.Rows.Item(22).AutoFit()

$ws.Range("A23").Clear()

$ws.Range("B23").Value = "LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3110 -  Projeto Integrado em Engenharia de Materiais III  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

$excel.CutCopyMode = $false
